# a69_f43_bUPPachuca.xlsx -- "cambio de fracciones e historico"
# Update the quarterly report: new reporting period (Q3 2022), refreshed
# validation/update dates, corrected responsible-area name, a couple of
# worksheet formatting tweaks (row height / column widths) and refreshed
# cell selections left behind by the editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Reporte de Formatos"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Reporte de Formatos")
$ws1.Activate() | Out-Null

# Corrected name of the responsible area (was misspelled / outdated)
$ws1.Range("G8").Value = "Subdirección de Planeación y Presupuesto (UPP)"

# New reporting quarter: Jul 1 - Sep 30, 2022 (previously Apr 1 - Jun 30)
$ws1.Range("B8").Value = "2022-07-01"
$ws1.Range("C8").Value = "2022-09-30"

# Refreshed validation / update dates
$ws1.Range("H8").Value = "2022-10-10"
$ws1.Range("I8").Value = "2022-10-10"

# J8 (Nota) gets the same left-aligned / bordered formatting as the rest
# of the data row
$ws1.Range("J8").HorizontalAlignment = -4131

# Give row 3 (merged headers) a little more room
$ws1.Rows("3").RowHeight = 36

# Leave the selection where the editor left it
$ws1.Range("B16").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Tabla_397514"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Tabla_397514")
$ws2.Activate() | Out-Null
$ws2.Range("C14").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Tabla_397515"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Tabla_397515")
$ws3.Activate() | Out-Null
$ws3.Columns("B").ColumnWidth = 19.7
$ws3.Range("E11").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Tabla_397516"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Tabla_397516")
$ws4.Activate() | Out-Null
$ws4.Columns("B").ColumnWidth = 20.8
$ws4.Range("D9").Select() | Out-Null

# ---------------------------------------------------------------------
# Re-activate the first sheet/tab, matching the saved workbook state
# ---------------------------------------------------------------------
$ws1.Activate() | Out-Null
